$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Productos")

# Clear the leftover empty placeholder cells on row 17 so they
# no longer exist in the sheet (F17, G17, L17).
$ws.Range("F17").ClearContents()
$ws.Range("G17").ClearContents()
$ws.Range("L17").ClearContents()

# Add new row 18 with the new order data.
$ws.Range("A18").Value = 1971
$ws.Range("B18").Value = "Ignacio Rodriguez"
$ws.Range("C18").Value = "Estructura coplanar NOVOTEGRA"
$ws.Range("D18").Value = "MODULO FV JA SOLAR 535WP BLACK FRAME BIFACIAL 120 CELDAS"
$ws.Range("E18").Value = "6"
$ws.Range("F18").Value = ""
$ws.Range("G18").Value = ""
$ws.Range("H18").Value = "GOODWE GW3600-ES-20 híbrido monofásico"
$ws.Range("I18").Value = "1"
$ws.Range("J18").Value = ""
$ws.Range("K18").Value = ""
$ws.Range("L18").Value = ""
$ws.Range("M18").Value = "Sí"
$ws.Range("N18").Value = "2025-09-25T07:50:43.054Z"
